$wb = $excel.ActiveWorkbook

# Add a new row of data to the "Users" sheet: F00474 / 074 (a commission
# case used for "consultas de bonificacion"), matching the existing
# style/number-format of the rows above it (column C is a text-formatted
# numeric code, so "074" must stay text rather than becoming 74).
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A5").Value = "F00474"
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "074"

# Move the active tab from "Modulos" to "Users", and leave the selection
# on C7 there.
$ws.Activate()
$ws.Range("C7").Select()

$wb.Save()
